# SYDATA-Work order testcases: roll the "RMA Details Maintenance Grid" sheet
# forward from the RMA-J098 batch to the new RMA-MG41 batch.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2
$ws.Range("E2").Value = "RMA-MG41-001"
$ws.Range("F2").Value = "RMA-MG41-1-1"
$ws.Range("J2").Value = "a7s5f000000xL33AAE"

# Row 3
$ws.Range("E3").Value = "RMA-MG41-002"
$ws.Range("F3").Value = "RMA-MG41-1-2"
$ws.Range("J3").Value = "a7s5f000000xL34AAE"

# Row 4
$ws.Range("E4").Value = "RMA-MG41-003"
$ws.Range("F4").Value = "RMA-MG41-1-3"
$ws.Range("J4").Value = "a7s5f000000xL35AAE"

# Columns resized (best-fit) to accommodate the new values
$ws.Columns.Item(5).ColumnWidth = 14.65
$ws.Columns.Item(6).ColumnWidth = 14.3
$ws.Columns.Item(10).ColumnWidth = 19.5
